$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Linked List")

# New rows of LeetCode tracker data to append (rows 4-12), matching the
# "did linked list cycle detection leetcode" commit.
# Row 4 keeps a date-like string ("08/12/2025"); force it to stay text
# (not auto-converted to a date serial number) like the existing A2/A3 cells.
$ws.Range("A4").NumberFormat = "@"
$ws.Range("A4").Value = "08/12/2025"

$ws.Range("B4").Value = "Linked List Cycle Detection"
$ws.Range("C4").Value = "Linked List"
$ws.Range("D4").Value = "Easy"
$ws.Range("E4").Value = "No (needed hint but was able to do code)"
$ws.Range("F4").Value = "Yes"
$ws.Range("G4").Value = "Yes"

$ws.Range("B5").Value = "Reorder List"
$ws.Range("C5").Value = "Linked List"
$ws.Range("D5").Value = "Medium"
$ws.Range("F5").Value = "Yes"

$ws.Range("B6").Value = "Remove Nth Node From End of List"
$ws.Range("C6").Value = "Linked List"
$ws.Range("D6").Value = "Medium"
$ws.Range("F6").Value = "Yes"

$ws.Range("B7").Value = "Copy List with Random Pointer"
$ws.Range("C7").Value = "Linked List"
$ws.Range("D7").Value = "Medium"
$ws.Range("F7").Value = "Yes"

$ws.Range("B8").Value = "Add Two Numbers"
$ws.Range("C8").Value = "Linked List"
$ws.Range("D8").Value = "Medium"
$ws.Range("F8").Value = "Yes"

$ws.Range("B9").Value = "Find The Duplicate Number"
$ws.Range("C9").Value = "Linked List"
$ws.Range("D9").Value = "Medium"
$ws.Range("F9").Value = "Yes"

$ws.Range("B10").Value = "LRU Cache"
$ws.Range("C10").Value = "Linked List"
$ws.Range("D10").Value = "Medium"
$ws.Range("F10").Value = "Yes"

$ws.Range("B11").Value = "Merge K Sorted Lists"
$ws.Range("C11").Value = "Linked List"
$ws.Range("D11").Value = "Hard"
$ws.Range("F11").Value = "Yes"

$ws.Range("B12").Value = "Reverse Nodes In K Group"
$ws.Range("C12").Value = "Linked List"
$ws.Range("D12").Value = "Hard"
$ws.Range("F12").Value = "Yes"
